$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append 10 new fixtures (rows 252-261, ids 250-259) to the bottom of the
# "Poland Ekstraklasa" results table, following the existing layout:
#   column A  -> id (bold / centered / bordered style, like the header rows)
#   column E  -> Date (custom date-time number format)
#   all other columns are plain text / numbers, same as the rows above them
# ---------------------------------------------------------------------------

# Carry the existing formatting of column A (bold/border "id" style) and
# column E (date/time number format) down into the new rows before filling
# in the values, so the appended rows look exactly like the ones above them.
$ws.Range("A251").Copy()
$ws.Range("A252:A261").PasteSpecial(-4122)
$ws.Range("E251").Copy()
$ws.Range("E252:E261").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 252 (id 250) ---
$ws.Range("A252").Value2 = 250
$ws.Range("B252").Value2 = 6775590
$ws.Range("C252").Value2 = "Poland Ekstraklasa"
$ws.Range("D252").Value2 = "Poland Ekstraklasa"
$ws.Range("E252").Value2 = 45397.58333333334
$ws.Range("F252").Value2 = "Piast Gliwice"
$ws.Range("G252").Value2 = "Zaglebie Lubin"
$ws.Range("H252").Value2 = 2
$ws.Range("I252").Value2 = 0
$ws.Range("J252").Value2 = "H"
$ws.Range("K252").Value2 = 2.15
$ws.Range("L252").Value2 = 3.2
$ws.Range("M252").Value2 = 3.4
$ws.Range("N252").Value2 = 2.2
$ws.Range("O252").Value2 = 2.9
$ws.Range("P252").Value2 = 3.6
$ws.Range("Q252").Value2 = -0.25
$ws.Range("R252").Value2 = 1.925
$ws.Range("S252").Value2 = 1.925
$ws.Range("T252").Value2 = 2
$ws.Range("U252").Value2 = 1.9
$ws.Range("V252").Value2 = 1.95
$ws.Range("W252").Value2 = 1.2
$ws.Range("X252").Value2 = -1
$ws.Range("Y252").Value2 = -1
$ws.Range("Z252").Value2 = 0.925
$ws.Range("AA252").Value2 = -1
$ws.Range("AB252").Value2 = 0
$ws.Range("AC252").Value2 = -0.0

# --- Row 253 (id 251) ---
$ws.Range("A253").Value2 = 251
$ws.Range("B253").Value2 = 6775595
$ws.Range("C253").Value2 = "Poland Ekstraklasa"
$ws.Range("D253").Value2 = "Poland Ekstraklasa"
$ws.Range("E253").Value2 = 45401.54166666666
$ws.Range("F253").Value2 = "Korona Kielce"
$ws.Range("G253").Value2 = "Radomiak Radom"
$ws.Range("K253").Value2 = 2.375
$ws.Range("L253").Value2 = 3.2
$ws.Range("M253").Value2 = 3.1
$ws.Range("N253").Value2 = 2.375
$ws.Range("O253").Value2 = 3.2
$ws.Range("P253").Value2 = 3.1
$ws.Range("Q253").Value2 = -0.25
$ws.Range("R253").Value2 = 2
$ws.Range("S253").Value2 = 1.85
$ws.Range("T253").Value2 = 2.25
$ws.Range("U253").Value2 = 2.025
$ws.Range("V253").Value2 = 1.825
$ws.Range("W253").Value2 = 0
$ws.Range("X253").Value2 = 0
$ws.Range("Y253").Value2 = 0
$ws.Range("Z253").Value2 = 0
$ws.Range("AA253").Value2 = 0

# --- Row 254 (id 252) ---
$ws.Range("A254").Value2 = 252
$ws.Range("B254").Value2 = 6775598
$ws.Range("C254").Value2 = "Poland Ekstraklasa"
$ws.Range("D254").Value2 = "Poland Ekstraklasa"
$ws.Range("E254").Value2 = 45401.64583333334
$ws.Range("F254").Value2 = "Rakow Czestochowa"
$ws.Range("G254").Value2 = "Gornik Zabrze"
$ws.Range("K254").Value2 = 1.571
$ws.Range("L254").Value2 = 4
$ws.Range("M254").Value2 = 5.75
$ws.Range("N254").Value2 = 1.571
$ws.Range("O254").Value2 = 4
$ws.Range("P254").Value2 = 5.75
$ws.Range("Q254").Value2 = -1
$ws.Range("R254").Value2 = 2.025
$ws.Range("S254").Value2 = 1.825
$ws.Range("T254").Value2 = 2.5
$ws.Range("U254").Value2 = 1.925
$ws.Range("V254").Value2 = 1.925
$ws.Range("W254").Value2 = 0
$ws.Range("X254").Value2 = 0
$ws.Range("Y254").Value2 = 0
$ws.Range("Z254").Value2 = 0
$ws.Range("AA254").Value2 = 0

# --- Row 255 (id 253) ---
$ws.Range("A255").Value2 = 253
$ws.Range("B255").Value2 = 6775597
$ws.Range("C255").Value2 = "Poland Ekstraklasa"
$ws.Range("D255").Value2 = "Poland Ekstraklasa"
$ws.Range("E255").Value2 = 45402.41666666666
$ws.Range("F255").Value2 = "Pogon Szczecin"
$ws.Range("G255").Value2 = "Piast Gliwice"
$ws.Range("K255").Value2 = 1.85
$ws.Range("L255").Value2 = 3.5
$ws.Range("M255").Value2 = 4.2
$ws.Range("N255").Value2 = 1.85
$ws.Range("O255").Value2 = 3.5
$ws.Range("P255").Value2 = 4.2
$ws.Range("Q255").Value2 = -0.5
$ws.Range("R255").Value2 = 1.875
$ws.Range("S255").Value2 = 1.975
$ws.Range("T255").Value2 = 2.25
$ws.Range("U255").Value2 = 1.825
$ws.Range("V255").Value2 = 2.025
$ws.Range("W255").Value2 = 0
$ws.Range("X255").Value2 = 0
$ws.Range("Y255").Value2 = 0
$ws.Range("Z255").Value2 = 0
$ws.Range("AA255").Value2 = 0

# --- Row 256 (id 254) ---
$ws.Range("A256").Value2 = 254
$ws.Range("B256").Value2 = 6774472
$ws.Range("C256").Value2 = "Poland Ekstraklasa"
$ws.Range("D256").Value2 = "Poland Ekstraklasa"
$ws.Range("E256").Value2 = 45402.52083333334
$ws.Range("F256").Value2 = "Ruch Chorzow"
$ws.Range("G256").Value2 = "Widzew Lodz"
$ws.Range("K256").Value2 = 2.55
$ws.Range("L256").Value2 = 3.3
$ws.Range("M256").Value2 = 2.7
$ws.Range("N256").Value2 = 2.55
$ws.Range("O256").Value2 = 3.3
$ws.Range("P256").Value2 = 2.7
$ws.Range("Q256").Value2 = 0
$ws.Range("R256").Value2 = 1.85
$ws.Range("S256").Value2 = 2
$ws.Range("T256").Value2 = 2.5
$ws.Range("U256").Value2 = 2
$ws.Range("V256").Value2 = 1.85
$ws.Range("W256").Value2 = 0
$ws.Range("X256").Value2 = 0
$ws.Range("Y256").Value2 = 0
$ws.Range("Z256").Value2 = 0
$ws.Range("AA256").Value2 = 0

# --- Row 257 (id 255) ---
$ws.Range("A257").Value2 = 255
$ws.Range("B257").Value2 = 6775594
$ws.Range("C257").Value2 = "Poland Ekstraklasa"
$ws.Range("D257").Value2 = "Poland Ekstraklasa"
$ws.Range("E257").Value2 = 45402.625
$ws.Range("F257").Value2 = "Zaglebie Lubin"
$ws.Range("G257").Value2 = "Jagiellonia Bialystok"
$ws.Range("K257").Value2 = 2.875
$ws.Range("L257").Value2 = 3.4
$ws.Range("M257").Value2 = 2.375
$ws.Range("N257").Value2 = 2.875
$ws.Range("O257").Value2 = 3.4
$ws.Range("P257").Value2 = 2.375
$ws.Range("Q257").Value2 = 0.25
$ws.Range("R257").Value2 = 1.775
$ws.Range("S257").Value2 = 2.1
$ws.Range("T257").Value2 = 2.75
$ws.Range("U257").Value2 = 2
$ws.Range("V257").Value2 = 1.85
$ws.Range("W257").Value2 = 0
$ws.Range("X257").Value2 = 0
$ws.Range("Y257").Value2 = 0
$ws.Range("Z257").Value2 = 0
$ws.Range("AA257").Value2 = 0

# --- Row 258 (id 256) ---
$ws.Range("A258").Value2 = 256
$ws.Range("B258").Value2 = 6850054
$ws.Range("C258").Value2 = "Poland Ekstraklasa"
$ws.Range("D258").Value2 = "Poland Ekstraklasa"
$ws.Range("E258").Value2 = 45403.3125
$ws.Range("F258").Value2 = "Cracovia Krakow"
$ws.Range("G258").Value2 = "Puszcza Niepolomice"
$ws.Range("K258").Value2 = 1.7
$ws.Range("L258").Value2 = 3.8
$ws.Range("M258").Value2 = 4.75
$ws.Range("N258").Value2 = 1.7
$ws.Range("O258").Value2 = 3.8
$ws.Range("P258").Value2 = 4.75
$ws.Range("Q258").Value2 = -0.75
$ws.Range("R258").Value2 = 1.925
$ws.Range("S258").Value2 = 1.925
$ws.Range("T258").Value2 = 2.5
$ws.Range("U258").Value2 = 1.975
$ws.Range("V258").Value2 = 1.875
$ws.Range("W258").Value2 = 0
$ws.Range("X258").Value2 = 0
$ws.Range("Y258").Value2 = 0
$ws.Range("Z258").Value2 = 0
$ws.Range("AA258").Value2 = 0

# --- Row 259 (id 257) ---
$ws.Range("A259").Value2 = 257
$ws.Range("B259").Value2 = 6830603
$ws.Range("C259").Value2 = "Poland Ekstraklasa"
$ws.Range("D259").Value2 = "Poland Ekstraklasa"
$ws.Range("E259").Value2 = 45403.41666666666
$ws.Range("F259").Value2 = "LKS Lodz"
$ws.Range("G259").Value2 = "Lech Poznan"
$ws.Range("K259").Value2 = 4.75
$ws.Range("L259").Value2 = 4
$ws.Range("M259").Value2 = 1.65
$ws.Range("N259").Value2 = 4.75
$ws.Range("O259").Value2 = 4
$ws.Range("P259").Value2 = 1.65
$ws.Range("Q259").Value2 = 0.75
$ws.Range("R259").Value2 = 2
$ws.Range("S259").Value2 = 1.85
$ws.Range("T259").Value2 = 2.75
$ws.Range("U259").Value2 = 1.975
$ws.Range("V259").Value2 = 1.875
$ws.Range("W259").Value2 = 0
$ws.Range("X259").Value2 = 0
$ws.Range("Y259").Value2 = 0
$ws.Range("Z259").Value2 = 0
$ws.Range("AA259").Value2 = 0

# --- Row 260 (id 258) ---
$ws.Range("A260").Value2 = 258
$ws.Range("B260").Value2 = 6775596
$ws.Range("C260").Value2 = "Poland Ekstraklasa"
$ws.Range("D260").Value2 = "Poland Ekstraklasa"
$ws.Range("E260").Value2 = 45403.52083333334
$ws.Range("F260").Value2 = "Legia Warsaw"
$ws.Range("G260").Value2 = "Slask Wroclaw"
$ws.Range("K260").Value2 = 1.7
$ws.Range("L260").Value2 = 3.6
$ws.Range("M260").Value2 = 5
$ws.Range("N260").Value2 = 1.7
$ws.Range("O260").Value2 = 3.6
$ws.Range("P260").Value2 = 5
$ws.Range("Q260").Value2 = -0.75
$ws.Range("R260").Value2 = 1.925
$ws.Range("S260").Value2 = 1.925
$ws.Range("T260").Value2 = 2.5
$ws.Range("U260").Value2 = 2.05
$ws.Range("V260").Value2 = 1.8
$ws.Range("W260").Value2 = 0
$ws.Range("X260").Value2 = 0
$ws.Range("Y260").Value2 = 0
$ws.Range("Z260").Value2 = 0
$ws.Range("AA260").Value2 = 0

# --- Row 261 (id 259) ---
$ws.Range("A261").Value2 = 259
$ws.Range("B261").Value2 = 6885526
$ws.Range("C261").Value2 = "Poland Ekstraklasa"
$ws.Range("D261").Value2 = "Poland Ekstraklasa"
$ws.Range("E261").Value2 = 45404.58333333334
$ws.Range("F261").Value2 = "Warta Poznan"
$ws.Range("G261").Value2 = "Stal Mielec"
$ws.Range("K261").Value2 = 2.1
$ws.Range("L261").Value2 = 3.1
$ws.Range("M261").Value2 = 3.8
$ws.Range("N261").Value2 = 2.1
$ws.Range("O261").Value2 = 3.1
$ws.Range("P261").Value2 = 3.8
$ws.Range("Q261").Value2 = -0.25
$ws.Range("R261").Value2 = 1.8
$ws.Range("S261").Value2 = 2.05
$ws.Range("T261").Value2 = 2
$ws.Range("U261").Value2 = 2.025
$ws.Range("V261").Value2 = 1.825
$ws.Range("W261").Value2 = 0
$ws.Range("X261").Value2 = 0
$ws.Range("Y261").Value2 = 0
$ws.Range("Z261").Value2 = 0
$ws.Range("AA261").Value2 = 0

